# Add season-record columns (Wins / Losses / Ties) to the roster table.
#
# The sheet currently ends at column AC (dimension A1:AC49). We append
# three new columns - AD "Wins", AE "Losses", AF "Ties" - with a header
# in row 1 and a value for every data row (2-49).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Reuse the existing bold/centered/bordered header formatting (the same
# style already applied to A1:AC1) by copying it onto the new headers.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null

# --- Data rows ---------------------------------------------------------
# Every roster row gets the team's 1998 season record: 90 wins, 73
# losses, 0 ties.
for ($r = 2; $r -le 49; $r++) {
    $ws.Cells.Item($r, 30).Value = 90   # AD - Wins
    $ws.Cells.Item($r, 31).Value = 73   # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF - Ties
}
